$d = $word.ActiveDocument

# --- Part 1: remove the stray _GoBack bookmark after the first <w:br/> ---
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# --- Part 2: append the new "Discrepancies" content at the end of the body ---
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $last.Range

# Paragraph: "Discrepancies:"
$r.InsertAfter("Discrepancies:")
$r.InsertParagraphAfter()
$r.Collapse(0)

# Paragraph: "If you tagged something and it is not showing up on the dashboard, troubleshooting steps include"
$r.InsertAfter("If you tagged something and it is not showing up on the dashboard, troubleshooting steps include")
$r.InsertParagraphAfter()
$r.Collapse(0)

# Paragraph: bulleted list item (ListParagraph style, same numbering as elsewhere in the doc)
$r.InsertAfter("Making sure that there is a PD Person account for that email ID AND that person account is linked to the deal")

$bulletPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$bulletPara.Style = "List Paragraph"
$bulletPara.Range.ListFormat.ApplyListTemplateWithLevel($d.ListTemplates.Item(1))

Write-Output "Done"
